# Extend testing data: add a new logbook row (row 11) and normalise the
# B9 driver-name cell back to the plain/default font.
#
# Target (per the OOXML diff):
#   - New shared string "Bob" (driver name)
#   - New row 11: AurXX1 / Bob / 2001-03-05 / 16:00 / 13.2 / 2001-03-05 / 17:14
#   - B9 loses its one-off font override (back to the same look as the rest
#     of column B)
#   - Date/time number-format codes become upper-case (DD/MM/YY, HH:MM)
#   - Final selection ends up on G12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new logbook entry ------------------------------------------
$ws.Cells.Item(11, 1).Value = "AurXX1"
$ws.Cells.Item(11, 2).Value = "Bob"
$ws.Cells.Item(11, 3).Value = 36952
$ws.Cells.Item(11, 4).Value = 0.666666666666667
$ws.Cells.Item(11, 5).Value = 13.2
$ws.Cells.Item(11, 6).Value = 36952
$ws.Cells.Item(11, 7).Value = 0.718055555555556

# --- B9: drop the special one-off font, back to the normal column font --
$ws.Cells.Item(9, 2).Font.Name = "Arial"
$ws.Cells.Item(9, 2).Font.Size = 10

# --- Date / time display formats, upper-cased ----------------------------
$ws.Range("C2:C11").NumberFormat = "DD/MM/YY"
$ws.Range("F2:F11").NumberFormat = "DD/MM/YY"
$ws.Range("D2:D11").NumberFormat = "HH:MM"
$ws.Range("G2:G11").NumberFormat = "HH:MM"

# --- Leave the selection where the editor's last action left it ---------
$ws.Range("G12").Select() | Out-Null
